$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.159.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.232.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.82%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.39%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.224.02'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.83%  '

$ws.Range("E9").Value = '  +4.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.07'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.21%  '

$ws.Range("E11").Value = '  +5.62%  '

$ws.Range("E12").Value = '  +5.74%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.42'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.64%  '

$ws.Range("E14").Value = '  +5.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.757.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.257.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.49%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '542.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +11.57%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.234.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.86%  '

$ws.Range("E19").Value = '  +3.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.64'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.43%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.745'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.80'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.36%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.43'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +21.59%  '

$ws.Range("E28").Value = '  +8.39%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.27'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.84%  '

$ws.Range("E31").Value = '  +6.64%  '

$ws.Range("E32").Value = '  -0.24%  '

$ws.Range("E33").Value = '  +5.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '565.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.67%  '

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.38'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.78%  '

$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.65'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.46%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0460'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '54.73'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.51%  '

$ws.Range("E39").Value = '  +7.68%  '

$ws.Range("E40").Value = '  +7.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.212.14'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +11.08%  '

$ws.Range("E42").Value = '  +4.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.62'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.67%  '

$ws.Range("E44").Value = '  +17.19%  '

$ws.Range("E45").Value = '  +13.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.999'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("E48").Value = '  +3.89%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.93'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.27%  '

$ws.Range("E50").Value = '  +4.00%  '

$ws.Range("E51").Value = '  +8.23%  '
